$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.874.99'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '3.330.06'
$ws.Range("E3").Value = '  -3.07%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '231.24'
$ws.Range("E5").Value = '  -2.63%  '
$ws.Range("D6").Value = '615.23'
$ws.Range("E6").Value = '  -3.76%  '
$ws.Range("E7").Value = '  -1.99%  '
$ws.Range("D8").Value = '0.384'
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.951'
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("D11").Value = '3.327.11'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").Value = '42.57'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = '0.194'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '6.01'
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").Value = '91.701.86'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").Value = '3.953.84'
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("E17").Value = '  -3.67%  '
$ws.Range("D18").Value = '8.06'
$ws.Range("E18").Value = '  -3.78%  '
$ws.Range("D19").Value = '3.325.36'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("D20").Value = '17.25'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").Value = '10.88'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("D22").Value = '3.41'
$ws.Range("E22").Value = '  +7.69%  '
$ws.Range("D23").Value = '493.14'
$ws.Range("E23").Value = '  -0.81%  '
$ws.Range("D24").Value = '0.448'
$ws.Range("E24").Value = '  -10.30%  '
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("D26").Value = '6.22'
$ws.Range("E26").Value = '  -5.66%  '
$ws.Range("D27").Value = '91.47'
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("D28").Value = '11.84'
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("D29").Value = '3.504.62'
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").Value = '11.04'
$ws.Range("E31").Value = '  -6.10%  '
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("E33").Value = '  -4.85%  '
$ws.Range("D34").Value = '0.992'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("E35").Value = '  -5.29%  '
$ws.Range("D36").Value = '28.13'
$ws.Range("E36").Value = '  -6.91%  '
$ws.Range("D37").Value = '0.526'
$ws.Range("E37").Value = '  -5.73%  '
$ws.Range("D38").Value = '561.27'
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("E42").Value = '  -5.54%  '
$ws.Range("D43").Value = '0.864'
$ws.Range("E43").Value = '  -7.42%  '
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0414'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D46").Value = '3.59'
$ws.Range("E46").Value = '  +2.81%  '
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").Value = '1.66'
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").Value = '5.39'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D49").Value = "'2.10"
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '7.91'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '51.33'
$ws.Range("E51").Value = '  -3.21%  '
